# CC_CallCenterData.xlsx — "Call center and BE updated"
#
# 1) Sheet1 (NormalSingleRoomBookingData) gains a bunch of "Modify booking
#    confirmation" columns (T:AF) plus a couple of swapped room-name values,
#    and two stray text cells on row 3 (Y3/Z3).
# 2) A brand-new sheet "ModifyBookingConfirmationPage" is appended after it,
#    holding a small adults/children/guest-name table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1: new header cells T1:AF1 (existing A1:S1 headers are untouched)
# ---------------------------------------------------------------------
$ws1.Range("T1").Value = "ConfirmationCode"
$ws1.Range("U1").Value = "RoomBefore"
$ws1.Range("V1").Value = "RoomAfter"
$ws1.Range("W1").Value = "DatesBefore"
$ws1.Range("X1").Value = "DatesAfter"
$ws1.Range("Y1").Value = "NoOfAdultAndChildBefore"
$ws1.Range("Z1").Value = "NoOfAdultAndChildAfter"
$ws1.Range("AA1").Value = "GuestDetailsBefore"
$ws1.Range("AB1").Value = "GuestDetailsAfter"
$ws1.Range("AC1").Value = "BedTypeBefore"
$ws1.Range("AD1").Value = "BedTypeAfter"
$ws1.Range("AE1").Value = "OtherServicesBefore"
$ws1.Range("AF1").Value = "OtherServicesAfter"

# Room names on row 2 swap: B2 was "Classic room" -> now "Standard Room";
# Q2 was "Standard room" -> now "Classic room".
$ws1.Range("B2").Value = "Standard Room"
$ws1.Range("Q2").Value = "Classic room"

# New data cells T2:AF2
$ws1.Range("T2").Value = "ALH100002653"
$ws1.Range("U2").Value = "Standard Room"
$ws1.Range("V2").Value = "Classic room"
$ws1.Range("W2").Value = "Thu, 2022-02-10➝Fri, 2022-02-11"
$ws1.Range("X2").Value = "Thu, 2022-02-10➝Sat, 2022-02-12"

# Y2/Z2 look like numbers ("2"/"3") but must stay plain shared-string text
# with no quote-prefix style, so type them with a leading apostrophe and
# then strip the resulting "quote prefix" style back to Normal.
$ws1.Range("Y2").Value = "'2"
$ws1.Range("Y2").Style = "Normal"
$ws1.Range("Z2").Value = "'3"
$ws1.Range("Z2").Style = "Normal"

$ws1.Range("AA2").Value = "rudraksh aggarwal"
$ws1.Range("AB2").Value = "john smith"

# AC2 is an empty-but-shared-string cell (t="s" pointing at an empty <t/>).
$ws1.Range("AC2").Value = "'"
$ws1.Range("AC2").Style = "Normal"

$ws1.Range("AD2").Value = "Classic Bed type"

# AE2/AF2 look like currency amounts; type them as text (leading apostrophe)
# and strip the resulting quote-prefix style back to Normal so they land as
# plain shared-string cells with no style override, same as the source file.
$ws1.Range("AE2").Value = "'$0.00"
$ws1.Range("AE2").Style = "Normal"
$ws1.Range("AF2").Value = "'$10.00"
$ws1.Range("AF2").Style = "Normal"

# Row 3 gains two more number-look-alike text cells (no quote-prefix style).
$ws1.Range("Y3").Value = "'0"
$ws1.Range("Y3").Style = "Normal"
$ws1.Range("Z3").Value = "'1"
$ws1.Range("Z3").Style = "Normal"

# ---------------------------------------------------------------------
# New sheet: ModifyBookingConfirmationPage (placed after sheet 1)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ModifyBookingConfirmationPage"

$ws2.Range("A1").Value = "Adults"
$ws2.Range("B1").Value = "Children"
$ws2.Range("C1").Value = "Child"
$ws2.Range("D1").Value = "FirstName"
$ws2.Range("E1").Value = "LastName"

$ws2.Range("A2").Value = "3 adults"
$ws2.Range("B2").Value = "'1"
$ws2.Range("C2").Value = "'1"
$ws2.Range("D2").Value = "john"
$ws2.Range("E2").Value = "smith"

$ws2.Range("F1").Select() | Out-Null

# Sheet1 stays the active tab; dimension grows to A1:AF5 and the selection
# moves from Q15 to G13 with no frozen/scrolled topLeftCell override.
$ws1.Activate() | Out-Null
$ws1.Range("G13").Select() | Out-Null
